# Scheduled GitHub Actions refresh of the crypto price table.
# Most rows just get updated Price/Volume(1h) text. In addition a new
# coin (BabyDogeCoin) now appears ahead of Cronos, so the former Cronos/
# Algorand/EnergySwap rows (49-51) shift down one slot each; EnergySwap,
# previously the last row of the A1:E51 table, drops off the bottom.
#
# Price/Volume cells are stored as plain text (t="inlineStr" in the
# original file) even when they look like a number (e.g. "0.998"), so
# for any such value we briefly mark the cell as Text before writing it
# and then clear that formatting again, to stop Excel's COM layer from
# auto-coercing the assignment into a numeric cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '26.518.36'
$ws.Range('E2').Value = '  +0.45%  '
$ws.Range('D3').Value = '1.630.81'
$ws.Range('E3').Value = '  +1.00%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.997'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '213.24'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.34%  '
$ws.Range('E6').Value = '  -0.18%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.998'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.21%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.247'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +0.23%  '
$ws.Range('E9').Value = '  +0.18%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.25'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +0.26%  '
$ws.Range('E11').Value = '  -0.55%  '
$ws.Range('D12').Value = '1.855.45'
$ws.Range('E12').Value = '  +0.71%  '
$ws.Range('D13').Value = '1.614.85'
$ws.Range('E13').Value = '  +0.04%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.04'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +0.04%  '
$ws.Range('E15').Value = '  +0.28%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.05'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -0.92%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '237.75'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +4.78%  '
$ws.Range('D18').Value = '26.510.24'
$ws.Range('E18').Value = '  +0.31%  '
$ws.Range('E19').Value = '  +4.46%  '
$ws.Range('D20').Value = '0.0₃0726'
$ws.Range('E20').Value = '  +0.00%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.998'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.05%  '
$ws.Range('E22').Value = '  -0.96%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.15'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.58%  '
$ws.Range('E24').Value = '  +2.61%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '147.09'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +1.33%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.11%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.10'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.90%  '
$ws.Range('E28').Value = '  +0.00%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.64'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +2.07%  '
$ws.Range('E30').Value = '  +0.20%  '
$ws.Range('E31').Value = '  -0.09%  '
$ws.Range('D32').Value = '1.524.75'
$ws.Range('E32').Value = '  +5.33%  '
$ws.Range('E33').Value = '  +1.25%  '
$ws.Range('E34').Value = '  -0.18%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.53'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +3.73%  '
$ws.Range('E36').Value = '  -0.08%  '
$ws.Range('E37').Value = '  +2.16%  '
$ws.Range('E38').Value = '  -0.02%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.838'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +0.20%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.89'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +0.82%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.999'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.09%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.21'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +0.76%  '
$ws.Range('D43').Value = '1.768.37'
$ws.Range('E43').Value = '  +0.84%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '63.12'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +2.23%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.762'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +0.22%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.909'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -0.95%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '90.49'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +2.93%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.52'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +1.83%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.0₆0102'
$ws.Range('E49').Value = '  -6.83%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0502'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.22%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0967'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.15%  '
